# Applies the "Fixed update to excel issue" change:
#  1. Rename header B1 on "Weekly Quantity" sheet -> "Weekly_PO_Qty"
#  2. Rename header B1 on "Monthly Trend" sheet -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet at the end with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet header rename ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet header rename ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header formatting (bold, centered, thin border) from the
# "Weekly Quantity" header row so the new sheet matches the workbook's style.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

# Reuse the existing date-cell formatting for column A (rows 2-15).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45487.99999999999, 16, 8.421966947125298, 24.58239964447902),
    @(45543.99999999999, 22, 14.11813412641464, 29.62336000187957),
    @(45557.99999999999, 23, 14.91905367409839, 31.33968963808949),
    @(45571.99999999999, 25, 16.10451597024165, 33.16167655057733),
    @(45592.99999999999, 27, 19.28337430162422, 34.80659079209495),
    @(45599.99999999999, 27, 19.03452573820964, 35.55207548048907),
    @(45606.99999999999, 28, 19.91505432480849, 36.12349049136315),
    @(45613.99999999999, 29, 20.0854718630087, 37.37447419918724),
    @(45620.99999999999, 29, 20.79755291518212, 38.25132295941745),
    @(45627.99999999999, 30, 21.47203860971021, 38.48605364543577),
    @(45634.99999999999, 31, 22.17761155953454, 38.58776205562328),
    @(45641.99999999999, 31, 23.40852805773612, 39.7137320853988),
    @(45648.99999999999, 32, 24.10651371603088, 40.216323896104),
    @(45655.99999999999, 33, 24.52276486966306, 41.05795225955208)
)

$row = 2
foreach ($item in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $item[0]
    $wsForecast.Cells.Item($row, 2).Value = $item[1]
    $wsForecast.Cells.Item($row, 3).Value = $item[2]
    $wsForecast.Cells.Item($row, 4).Value = $item[3]

    $row = $row + 1
}

$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
